$p = $ppt.ActivePresentation

# Remove the slide titled "Rotation - issues" (sldId 278) from the deck.
# It is the 21st slide in the current slide order.
$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $title = ""
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $title = $shape.TextFrame.TextRange.Text
            break
        }
    }
    if ($title -eq "Rotation - issues") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $p.Slides.Item($target).Delete()
} else {
    $p.Slides.Item(21).Delete()
}
